$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The price-tracking sheet picked up a fresh scrape: the previous last row's
# timestamp got a (sub-microsecond) refresh, and a brand new row was appended
# with the latest reading for the same product/weight/price.
$ws.Range("A12").Value = 45812.39347450231

$ws.Range("A13").Value = 45813.39350918835
# Match the date/time number format already used by the rest of column A.
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat
$ws.Range("B13").Value = "EVOWHEY PROTEIN"
$ws.Range("C13").Value = "2Kg"
$ws.Range("D13").Value = "34,90€"
